$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source row (last existing data row) to replicate for the new rows below it
$src = $ws.Range("A101:J101")

# Copy the formatting/values of row 101 into the four new rows
$src.Copy($ws.Range("A102:J102"))
$src.Copy($ws.Range("A103:J103"))
$src.Copy($ws.Range("A104:J104"))
$src.Copy($ws.Range("A105:J105"))

# Update the date serial values in column A for the newly added rows
$ws.Range("A102").Value2 = 45658
$ws.Range("A103").Value2 = 45659
$ws.Range("A104").Value2 = 45660
$ws.Range("A105").Value2 = 45661
